$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 11
$ws.Range("B6").Value = "Vega Monumental Concepción"
$ws.Range("C6").Value = "Bíobío"
$ws.Range("D6").Value = 44714
$ws.Range("D6").NumberFormat = $ws.Range("D5").NumberFormat
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = "Otros"
$ws.Range("I6").Value = 100107001
$ws.Range("J6").Value = "Caqui"
$ws.Range("K6").Value = "Mankaki"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 14000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 14500
$ws.Range("Q6").Value = "$/caja 18 kilos granel"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 806
$ws.Range("T6").Value = 18
